# A new daily price record was collected for Mango at Vega Monumental
# Concepción. It is inserted as a new row at position 97, pushing the
# existing rows 97-201 down to 98-202 (the data block is not sorted by
# date, so the new record lands at the top of the range rather than at
# the very end).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 97:201 down to 98:202, leaving a blank row 97 (Excel copies
# the formatting, e.g. the date number format on column D, from the row
# above).
$ws.Rows("97:97").Insert()

# Populate the new row 97 with the new record.
$ws.Range("A97").Value = 11
$ws.Range("B97").Value = 'Vega Monumental Concepción'
$ws.Range("C97").Value = 'Bíobío'
$ws.Range("D97").Value = 45225
$ws.Range("E97").Value = 8
$ws.Range("F97").Value = 'Fruta'
$ws.Range("G97").Value = 100108
$ws.Range("H97").Value = 'Tropicales y subtropicales'
$ws.Range("I97").Value = 100108002
$ws.Range("J97").Value = 'Mango'
$ws.Range("K97").Value = 'Sin especificar'
$ws.Range("L97").Value = 'Primera'
$ws.Range("M97").Value = 200
$ws.Range("N97").Value = 9000
$ws.Range("O97").Value = 9500
$ws.Range("P97").Value = 9250
$ws.Range("Q97").Value = '$/bandeja 4 kilos'
$ws.Range("R97").Value = 'Brasil'
$ws.Range("S97").Value = 2312
$ws.Range("T97").Value = 4
